# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Worksheets item 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 337
$ws1.Range("F5").Value  = 1727
$ws1.Range("F8").Value  = 5
$ws1.Range("F11").Value = 4893
$ws1.Range("F17").Value = 178
$ws1.Range("F21").Value = 3839
$ws1.Range("F22").Value = 706
$ws1.Range("F23").Value = 648
$ws1.Range("F31").Value = 572
$ws1.Range("F34").Value = 929
$ws1.Range("F35").Value = 2443

# --- Sheet "全部类型" (Worksheets item 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 337
$ws4.Range("F5").Value  = 1727
$ws4.Range("F8").Value  = 5
$ws4.Range("F11").Value = 4893
$ws4.Range("F17").Value = 178
$ws4.Range("F21").Value = 3839
$ws4.Range("F22").Value = 706
$ws4.Range("F23").Value = 648
$ws4.Range("F31").Value = 572
$ws4.Range("F35").Value = 929
$ws4.Range("F36").Value = 2443
